$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A width (~16.5703125, closest achievable via ColumnWidth rounding) ---
$ws.Columns.Item(1).ColumnWidth = 15.65

# --- Row 4: shrink the huge-font row height, mark it as a custom height ---
$ws.Rows.Item(4).RowHeight = 29.25

# --- Row 6: new row with the small font and a very small custom height ---
$ws.Range("A6").Value2 = "Row with very small font and very small height"
$ws.Range("A6").Font.Size = $ws.Range("A3").Font.Size
$ws.Rows.Item(6).RowHeight = 6.75

# --- Row 8: new row with the huge font, wrapped text and a numeric value ---
$ws.Range("A8").Value2 = 123
$ws.Range("A8").Font.Size = $ws.Range("A4").Font.Size
$ws.Range("A8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 15.75

# --- View state: zoom in and move the selection to the new last cell ---
$aw = $ws.Application.ActiveWindow
$aw.Zoom = 238
$ws.Range("A8").Select()
